$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 66 ("Supervision Condition" mapping row).
# Excel copies formatting from the row above (row 65: s13/customFormat), giving the
# new row's cells the same styles as the target (A66 style 8, B66 style 13, E66 style 2).
$ws.Rows.Item(66).Insert()

# Fill in the new row's content.
$ws.Range("B66").Value = "Supervision Condition"
$ws.Range("E66").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Detention[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityCaseAssociation/nc:Activity/@structures:ref]/j:SupervisionAugmentation/j:SupervisionCondition/nc:ActivityDescriptionText"

# The inserted row picks up a stray C66 cell from the copied formatting; the target
# row has no C column entry, so clear it back out.
$ws.Range("C66").Clear()

# Match the row height used for the new row in the target workbook.
$ws.Rows.Item(66).RowHeight = 60

# Leave the cursor on C66, matching the author's final selection.
$ws.Range("C66").Select()
